$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.025.25"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.36%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.580.74"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.05%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "549.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.88%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.02%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  +2.21%  "
$ws.Range("E9").Value = "  -0.80%  "
$ws.Range("E10").Value = "  -1.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.58"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.29%  "
$ws.Range("E12").Value = "  -0.98%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.038.46"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.62"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.01%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "61.920.58"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.31%  "
$ws.Range("E16").Value = "  +0.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.581.19"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.14%  "
$ws.Range("E18").Value = "  -3.56%  "
$ws.Range("E19").Value = "  -0.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "338.35"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.47%  "
$ws.Range("E21").Value = "  -4.63%  "
$ws.Range("E22").Value = "  +0.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.493"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.71"
$ws.Range("D24").Style = "Normal"
$ws.Range("E25").Value = "  -0.85%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.17"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0840"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.12%  "
$ws.Range("E30").Value = "  +1.64%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.89"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.98%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "162.62"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.88"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.99%  "
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "19.22"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.79%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.43"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.22%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.80"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.95%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "329.43"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.14%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.07"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.912"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.97%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.95"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.42%  "
$ws.Range("E42").Value = "  -1.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.09"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.29%  "
$ws.Range("E44").Value = "  -0.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.607"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.55%  "
$ws.Range("B46").Value = "WhiteBITCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.94"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.114.99"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.85%  "
$ws.Range("B48").Value = "Hedera"
$ws.Range("C48").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0548"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.84%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.56"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.64%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0966"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.69%  "
